$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update elt_geld values (column H) to fix the test data
$ws.Range("H4").Value = 1005
$ws.Range("H6").Value = 1206
$ws.Range("H13").Value = 603

# Update the active cell selection to H1
$ws.Range("H1").Select()
